$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample was recorded for 2026/02/21 (土) at hour 20 -> rank 36.
# Insert it right after the existing 2026/02/21 rows (row 828), pushing
# every following row down by one.
$ws.Rows("828:828").Insert()

# Force column A to be entered as literal text (matches the rest of the
# date column, which is stored as text, not a real date) and then drop
# the temporary format back to the sheet default so no extra style is
# left behind on the cell.
$ws.Range("A828").NumberFormat = "@"
$ws.Range("A828").Value = "2026/02/21"
$ws.Range("A828").Style = "Normal"

$ws.Range("B828").Value = "土"
$ws.Range("C828").Value = 20
$ws.Range("D828").Value = 36
